# Apply the "subject and mail body with variable values and minor errors solved" edit.
# The sheet originally held a list of 6 sample recipients (rows 2-7); the author trimmed it
# down to a single, corrected sample row and renamed the header labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the first five data rows (rows 2-6); the former row 7 (Vinay Savla's data) shifts up
# to become the new row 2, which is the only data row kept.
$ws.Range("A2:D6").EntireRow.Delete()

# Rename the header row to the new, cleaner labels.
$ws.Range("A1").Value = "NAME"
$ws.Range("B1").Value = "SALARY"
$ws.Range("C1").Value = "POSITION"

# Fix the name typo/casing for the remaining sample row.
$ws.Range("A2").Value = "Vinay Savla"

# The hyperlink collection still references the old (now deleted) rows; clear them out and
# re-create the single remaining hyperlink on D2, restoring the Hyperlink cell style.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:savlavinay022@gmail.com")
$ws.Range("D2").Style = "Hyperlink"

# Match the saved selection/active cell from the final workbook state.
[void]$ws.Range("C5").Select()
